# Revert ""Fuel" sheet update for both CH and SIN"
# This reverts the earlier commit that updated KBOB-sourced natural-gas /
# electricity / biogas figures on the FUELS and ELECTRICITY sheets of the
# CH LCA_infrastructure workbook back to their CEA-sourced values, and
# removes the short-lived "Biogas" rows / shared strings again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# FUELS sheet — restore previous Natural Gas / Electricity figures,
# clear out Oil/Coal/Wood placeholder values, and drop the Biogas row
# (and the now-trailing blank style row) that the reverted commit added.
# ---------------------------------------------------------------------
$wsFuels = $wb.Worksheets.Item("FUELS")

# Natural Gas (row 2)
$wsFuels.Range("C2").Value = 1.403
$wsFuels.Range("D2").Value = 0.1
$wsFuels.Range("F2").Value = "from CEA, costs in USD-2015,"

# Electricity / GRID (row 3)
$wsFuels.Range("C3").Value = 2.63
$wsFuels.Range("D3").Value = 0.0413
$wsFuels.Range("F3").Value = "from CEA, costs in USD-2015,"

# Oil, Coal, Wood (rows 5-7) lose their values/sources again
$wsFuels.Range("C5:F7").ClearContents()

# The Biogas data row (8) and the trailing blank formatting row (9) go away
$wsFuels.Rows("8:9").Delete()

# ---------------------------------------------------------------------
# ELECTRICITY sheet — same GRID figures roll back, reference note goes
# back to the generic "KBOB 2019" string instead of the CH-Verbrauchermix
# one (which is why that shared string becomes unused and disappears).
# ---------------------------------------------------------------------
$wsElec = $wb.Worksheets.Item("ELECTRICITY")
$wsElec.Range("E3").Value = 2.63
$wsElec.Range("F3").Value = 0.0413
$wsElec.Range("H3").Value = "KBOB 2019, costs in USD-2015"
# E5/F5 are formulas off of E3:E4/F3:F4, so they recalc automatically.

# ---------------------------------------------------------------------
# Restore the previous window/selection state: FUELS was the active tab
# and zoomed to 150%; after the revert HEATING is active at 100% and each
# sheet's cursor lands back on its pre-update cell.
# ---------------------------------------------------------------------
$wsDHW = $wb.Worksheets.Item("DHW")
$wsDHW.Range("E1:G1048576").Select() | Out-Null

$wsCooling = $wb.Worksheets.Item("COOLING")
$wsCooling.Range("E1:G1048576").Select() | Out-Null

$wsElec.Range("D27").Select() | Out-Null

$wsFuels.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 100
$wsFuels.Range("F18").Select() | Out-Null

$wsHeating = $wb.Worksheets.Item("HEATING")
$wsHeating.Activate() | Out-Null
$wsHeating.Range("B35").Select() | Out-Null

Write-Host "Revert applied"
